# Applies the "Added 'Date and Time' and 'Cycle_count' parameters" edit.
# The sheet is effectively re-laid-out: a new "Date and Time" row is
# inserted at the top, a new "Total SOC consumed(%)" row is inserted
# after "Mode", a new "Cycle Count of battery" row is inserted after the
# electricity-consumption row, the Cell-Voltage and cell-temperature rows
# are reordered, several labels/values are corrected, and two new
# "Time spent in ..." rows are appended at the bottom. Rather than doing
# piecemeal row-inserts, we just rewrite column A/B for every target row
# (1..45) with the final content, then fix up the one cell that needs the
# special time format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any pre-existing per-cell number formatting in the target area so
# stale formatting (e.g. the old B1 time format) doesn't leak onto the
# wrong row once everything shifts down.
$ws.Range("A1:B45").Style = "Normal"

$rows = @(
    @{ Row=1; A="Date and Time"; BKind="str"; BVal="2024-03-11 10:56:55.193000 to 2024-03-11 11:45:44.722000" },
    @{ Row=2; A="Total time taken for the ride"; BKind="num"; BVal="0.03393528935185185" },
    @{ Row=3; A="Actual Ampere-hours (Ah)"; BKind="num"; BVal="28.38003055555555" },
    @{ Row=4; A="Actual Watt-hours (Wh)"; BKind="num"; BVal="1430.000783129167" },
    @{ Row=5; A="Starting SoC (Ah)"; BKind="num"; BVal="38.841" },
    @{ Row=6; A="Ending SoC (Ah)"; BKind="num"; BVal="10.178" },
    @{ Row=7; A="Starting SoC (%)"; BKind="num"; BVal="98" },
    @{ Row=8; A="Ending SoC (%)"; BKind="num"; BVal="25" },
    @{ Row=9; A="Total distance covered (km)"; BKind="num"; BVal="31.96949405624171" },
    @{ Row=10; A="Total energy consumption(WH/KM)"; BKind="num"; BVal="44.73016622075612" },
    @{ Row=11; A="Total SOC consumed(%)"; BKind="num"; BVal="73" },
    @{ Row=12; A="Mode"; BKind="str"; BVal="Custom mode`n97.11%`nEco mode`n1.34%`nSports mode`n0.08%" },
    @{ Row=13; A="Peak Power(kW)"; BKind="num"; BVal="5437.23852" },
    @{ Row=14; A="Average Power(kW)"; BKind="num"; BVal="-1767.858111011332" },
    @{ Row=15; A="Total Energy Regenerated(kWh)"; BKind="num"; BVal="0.045046905" },
    @{ Row=16; A="Regenerative Effectiveness(%)"; BKind="num"; BVal="0.003150032261478203" },
    @{ Row=17; A="Highest Cell Voltage(V)"; BKind="num"; BVal="3.326" },
    @{ Row=18; A="Lowest Cell Voltage(V)"; BKind="num"; BVal="3" },
    @{ Row=19; A="Difference in Cell Voltage(V)"; BKind="num"; BVal="0.3260000000000001" },
    @{ Row=20; A="Minimum Temperature(C)"; BKind="num"; BVal="27" },
    @{ Row=21; A="Maximum Temperature(C)"; BKind="num"; BVal="41" },
    @{ Row=22; A="Difference in Temperature(C)"; BKind="num"; BVal="14" },
    @{ Row=23; A="Maximum Fet Temperature-BMS(C)"; BKind="num"; BVal="67" },
    @{ Row=24; A="Maximum Afe Temperature-BMS(C)"; BKind="num"; BVal="60" },
    @{ Row=25; A="Maximum PCB Temperature-BMS(C)"; BKind="num"; BVal="59" },
    @{ Row=26; A="Maximum MCU Temperature(C)"; BKind="num"; BVal="59" },
    @{ Row=27; A="Maximum Motor Temperature(C)"; BKind="num"; BVal="97" },
    @{ Row=28; A="Abnormal Motor Temperature Detected(C)"; BKind="num"; BVal="0" },
    @{ Row=29; A="highest cell temp(C)"; BKind="num"; BVal="41" },
    @{ Row=30; A="lowest cell temp(C)"; BKind="num"; BVal="27" },
    @{ Row=31; A="Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"; BKind="num"; BVal="14" },
    @{ Row=32; A="Battery Voltage(V)"; BKind="num"; BVal="53" },
    @{ Row=33; A="Total energy charged(kWh)"; BKind="num"; BVal="1.504141619444444" },
    @{ Row=34; A="Electricity consumption units(kW)"; BKind="num"; BVal="1.42502427187021e-07" },
    @{ Row=35; A="Cycle Count of battery"; BKind="num"; BVal="26" },
    @{ Row=36; A="Idling time percentage"; BKind="num"; BVal="20.11036174126303" },
    @{ Row=37; A="Time spent in 0-10 km/h"; BKind="num"; BVal="3.653478558805497" },
    @{ Row=38; A="Time spent in 10-20 km/h"; BKind="num"; BVal="5.106935478053883" },
    @{ Row=39; A="Time spent in 20-30 km/h"; BKind="num"; BVal="10.85945107656797" },
    @{ Row=40; A="Time spent in 30-40 km/h"; BKind="num"; BVal="10.52043134850507" },
    @{ Row=41; A="Time spent in 40-50 km/h"; BKind="num"; BVal="7.837126266815739" },
    @{ Row=42; A="Time spent in 50-60 km/h"; BKind="num"; BVal="10.09124679914884" },
    @{ Row=43; A="Time spent in 60-70 km/h"; BKind="num"; BVal="13.27226169437732" },
    @{ Row=44; A="Time spent in 70-80 km/h"; BKind="num"; BVal="18.37559057957947" },
    @{ Row=45; A="Time spent in 80-90 km/h"; BKind="num"; BVal="0" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    if ($item.BKind -eq "num") {
        $ws.Cells.Item($r, 2).Value = [double]$item.BVal
    } else {
        $ws.Cells.Item($r, 2).Value = $item.BVal
    }
}

# Row 2 ("Total time taken for the ride") keeps the original [hh]:mm:ss
# duration format that used to live on row 1.
$ws.Cells.Item(2, 2).NumberFormat = "[hh]:mm:ss"

Write-Host "Applied Date and Time / Cycle Count edit"
